$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.26"
$ws.Range("E2").Value = "'0.60%"
$ws.Range("D3").Value = "'27.15"
$ws.Range("E3").Value = "'1.98%"
$ws.Range("D4").Value = "'4.843"
$ws.Range("E4").Value = "'-1.02%"
$ws.Range("D5").Value = "'0.06397"
$ws.Range("D6").Value = "'6.927"
$ws.Range("E6").Value = "'0.21%"
$ws.Range("D7").Value = "'1.204"
$ws.Range("E7").Value = "'-5.19%"
$ws.Range("D8").Value = "'0.8790"
$ws.Range("E8").Value = "'0.10%"
$ws.Range("D9").Value = "'0.1515"
$ws.Range("E9").Value = "'4.24%"
$ws.Range("D10").Value = "'0.05025"
$ws.Range("E10").Value = "'0.79%"
$ws.Range("D11").Value = "'0.07549"
$ws.Range("E11").Value = "'3.07%"
$ws.Range("D12").Value = "'0.02971"
$ws.Range("E12").Value = "'-4.91%"
$ws.Range("D13").Value = "'0.09002"
$ws.Range("E13").Value = "'-0.57%"
$ws.Range("D14").Value = "'0.001573"
$ws.Range("E14").Value = "'-0.39%"
$ws.Range("E15").Value = "'1.69%"
$ws.Range("D16").Value = "'0.006198"
$ws.Range("E16").Value = "'2.28%"
$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'0.39%"
$ws.Range("D18").Value = "'3.312"
$ws.Range("E18").Value = "'-1.40%"
$ws.Range("E19").Value = "'0.55%"
$ws.Range("E20").Value = "'-0.96%"
$ws.Range("D21").Value = "'0.1360"
$ws.Range("E21").Value = "'2.38%"
$ws.Range("D22").Value = "'3.931"
$ws.Range("E22").Value = "'0.67%"
$ws.Range("D23").Value = "'0.04421"
$ws.Range("E23").Value = "'0.20%"
$ws.Range("E24").Value = "'-0.39%"
$ws.Range("D25").Value = "'0.003852"
$ws.Range("E25").Value = "'-12.58%"
$ws.Range("D26").Value = "'0.0001203"
$ws.Range("E26").Value = "'0.12%"
$ws.Range("D27").Value = "'0.0001941"
$ws.Range("E27").Value = "'13.81%"
$ws.Range("D40").Value = "'0.04133"
$ws.Range("E40").Value = "'2.44%"
$ws.Range("D41").Value = "'0.006830"
$ws.Range("E41").Value = "'2.64%"
$ws.Range("E42").Value = "'0.76%"
$ws.Range("D43").Value = "'0.002165"
$ws.Range("E43").Value = "'2.98%"
$ws.Range("D44").Value = "'0.01150"
$ws.Range("E44").Value = "'-2.56%"
$ws.Range("D45").Value = "'0.00005163"
$ws.Range("E45").Value = "'-3.36%"
$ws.Range("D46").Value = "'1.650"
$ws.Range("E46").Value = "'-38.15%"
$ws.Range("D47").Value = "'0.02005"
$ws.Range("E47").Value = "'-0.04%"
